$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 10:22"

# --- Reorder "Polonia" ahead of "Dinamarca" / "Chequia" in the countries list ---
# Row 29 was Dinamarca, now becomes Polonia (with refreshed stats)
$ws.Range("A29").Value = "Polonia"
$ws.Range("B29").Value = 5742
$ws.Range("C29").Value = 167
$ws.Range("D29").Value = 318
$ws.Range("E29").Value = 5249
$ws.Range("F29").Value = 160
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 175

# Row 30 was Chequia, now becomes Dinamarca
$ws.Range("A30").Value = "Dinamarca"
$ws.Range("B30").Value = 5635
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 1736
$ws.Range("E30").Value = 3662
$ws.Range("F30").Value = 120
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 237

# Row 31 was Polonia, now becomes Chequia
$ws.Range("A31").Value = "Chequia"
$ws.Range("B31").Value = 5589
$ws.Range("C31").Value = 20
$ws.Range("D31").Value = 309
$ws.Range("E31").Value = 5167
$ws.Range("F31").Value = 98
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 113

# --- Update other case statistics (row 19 Austria, row 34 Rumania, row 38 Filipinas) ---
$ws.Range("B19").Value = 13337
$ws.Range("C19").Value = 93
$ws.Range("E19").Value = 7802

$ws.Range("E34").Value = 4298
$ws.Range("G34").Value = 9
$ws.Range("H34").Value = 257

$ws.Range("B38").Value = 4195
$ws.Range("C38").Value = 119
$ws.Range("D38").Value = 140
$ws.Range("E38").Value = 3834
$ws.Range("G38").Value = 18
$ws.Range("H38").Value = 221
